$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last test case's "Resultado" cell (B9) should now report an HTTP 403
# response instead of 200, matching the "No Ingresa a la BD" scenario used
# elsewhere in the sheet.
$ws.Range("B9").Value = "403`nNo Ingresa a la  BD"

# Reflect the view state after the edit: B9 is now the active/selected
# cell, and the window is scrolled so row 6 is at the top.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B9").Select()
